$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# The "Rescatables" table (rows 6-10) needs a personal-data correction: the
# whole record (A:G) that used to sit in row 6 actually belongs at the
# bottom of the block (row 10), and the records that used to occupy rows
# 7-10 each move up one row. Net effect: rows 6-10 are cyclically rotated
# up by one row.

# Capture row 6's current values before anything shifts.
$a6 = $ws.Range("A6").Value()
$b6 = $ws.Range("B6").Value()
$c6 = $ws.Range("C6").Value()
$d6 = $ws.Range("D6").Value()
$e6 = $ws.Range("E6").Value()
$f6 = $ws.Range("F6").Value()
$g6 = $ws.Range("G6").Value()

# Deleting row 6 shifts the old rows 7-10 up to become rows 6-9.
$ws.Rows("6").Delete()

# Row 10 is now free; put the old row 6 record there.
$ws.Range("A10").Value = $a6
$ws.Range("B10").Value = $b6
$ws.Range("C10").Value = $c6
$ws.Range("D10").Value = $d6
$ws.Range("E10").Value = $e6
$ws.Range("F10").Value = $f6
$ws.Range("G10").Value = $g6
